$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new time-log entry on row 79 (values previously blank)
$ws.Range("A79").Value = 41920
$ws.Range("B79").Value = 0.5229166666666667
$ws.Range("C79").Value = 0.61458333333333337
$ws.Range("D79").Value = 15
$ws.Range("F79").Value = "Coding"

# The Delta column (E) uses a shared formula already spanning this row;
# make sure row 79 keeps the same (row-relative) formula as the rest of the column.
$ws.Range("E79").Formula = "=IF(AND(NOT(ISBLANK(B79)),NOT(ISBLANK(C79))), (C79-B79) * 24 - D79/60, """")"

$excel.CalculateFullRebuild()
$excel.Calculate()

# Refresh the pie chart on Sheet2 so its cached values pick up the new totals.
$ws2 = $wb.Worksheets.Item("Sheet2")
foreach ($co in $ws2.ChartObjects()) {
    $co.Chart.Refresh()
}

$ws.Select()
$ws.Range("C80").Select()

$wb.Save()
